$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection to I18
$ws.Range("I18").Select()

# Row 2: AE2:AH2 change from "K" to "C"
$ws.Range("AE2:AH2").Value = "C"

# Row 4: AE4:AH4 change from "C" to "K"
$ws.Range("AE4:AH4").Value = "K"
